$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row added: the "建议" (suggestions) prompt goes into the now-empty F9 cell first,
# then the existing F8 "任务清单" text is edited down to "任务".
$ws.Range("F9").Value = "建议/请精心提供20个最好的建议和具体细节。回答样式:{建议}-{建议的具体组成部分和细节内容}  Final output are in the following format:     - 段落 1     - 段落 2     - 段落 3"
$ws.Range("F8").Value = "任务/请精心设计20个主要任务和主要完成的步骤。回答样式:{任务介绍}-{详细介绍主要完成步骤和方法}  Final output are in the following format:     - 段落 1     - 段落 2     - 段落 3"

$ws.Range("F9").Select()
